$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two trailing monthly rows (41, 42) plus the blank placeholder row (43)
# are removed: row 41's data is cleared out (leaving an empty, still
# date-styled, placeholder row) and rows 42:43 are deleted outright so the
# sheet's used range shrinks from A1:E43 down to A1:E41.
$ws.Range("A41:D41").ClearContents()
$ws.Range("A42:E42").EntireRow.Delete()
$ws.Range("A42:E42").EntireRow.Delete()

# Reproduce the saved view state: scrolled down so row 32 is the first
# visible row, with A41:A42 selected (mirrors selecting the old rows
# 41:42 right before deleting them).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 32
$ws.Range("A41:A42").Select()
